{"js": "// Fix two small pt-BR translation/wording slips in the API documentation:\n//   1) \"A URL para acessar o API ...\"  ->  \"A URL para acessar a API ...\"\n//   2) \"Busca uma Sele\u00e7\u00e3o pelo nome ou palavra-chave\"\n//        -> \"Busca um Jogador pela sele\u00e7\u00e3o ou palavra-chave\"\n// (both occurrences are unique in the document, so an exact-text search is safe)\n\nconst body = context.document.body;\n\n// --- Edit 1: \"URL\" paragraph -------------------------------------------------\nconst urlMatches = body.search(\"A URL para acessar o API\", { matchCase: true });\nurlMatches.load(\"items\");\nawait context.sync();\n\nif (urlMatches.items.length > 0) {\n  urlMatches.items[0].insertText(\n    \"A URL para acessar a API\",\n    Word.InsertLocation.replace\n  );\n}\n\n// --- Edit 2: \"Jogadores\" (Busca ...) paragraph -------------------------------\nconst buscaMatches = body.search(\n  \"Busca uma Sele\u00e7\u00e3o pelo nome ou palavra-chave\",\n  { matchCase: true }\n);\nbuscaMatches.load(\"items\");\nawait context.sync();\n\nif (buscaMatches.items.length > 0) {\n  buscaMatches.items[0].insertText(\n    \"Busca um Jogador pela sele\u00e7\u00e3o ou palavra-chave\",\n    Word.InsertLocation.replace\n  );\n}\n\nawait context.sync();\n", "ps1": "# Fix two small pt-BR translation/wording slips in the API documentation:\n#   1) \"A URL para acessar o API ...\"  ->  \"A URL para acessar a API ...\"\n#   2) \"Busca uma Sele\u00e7\u00e3o pelo nome ou palavra-chave\"\n#        -> \"Busca um Jogador pela sele\u00e7\u00e3o ou palavra-chave\"\n# (both occurrences are unique in the document, so an exact-text Find is safe)\n\n$d = $word.ActiveDocument\n\n# --- Edit 1: \"URL\" paragraph -------------------------------------------------\n$range1 = $d.Content\n$find1 = $range1.Find\n$find1.Text = \"A URL para acessar o API\"\n$find1.MatchCase = $true\n$find1.MatchWholeWord = $false\nif ($find1.Execute()) {\n    $range1.Text = \"A URL para acessar a API\"\n}\n\n# --- Edit 2: \"Jogadores\" (Busca ...) paragraph -------------------------------\n$range2 = $d.Content\n$find2 = $range2.Find\n$find2.Text = \"Busca uma Sele\u00e7\u00e3o pelo nome ou palavra-chave\"\n$find2.MatchCase = $true\n$find2.MatchWholeWord = $false\nif ($find2.Execute()) {\n    $range2.Text = \"Busca um Jogador pela sele\u00e7\u00e3o ou palavra-chave\"\n}\n"}
